$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.555.83'
$ws.Range("E2").Value = '  -0.38%  '
$ws.Range("D3").Value = '1.634.88'
$ws.Range("E3").Value = '  +0.21%  '
$ws.Range("E4").Value = '  +0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.69'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.503'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.03%  '
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.250'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.47%  '
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0626'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.35%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.84'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.73%  '
$ws.Range("E11").Value = '  +0.37%  '
$ws.Range("D12").Value = '1.861.46'
$ws.Range("E12").Value = '  +0.12%  '
$ws.Range("D13").Value = '1.647.28'
$ws.Range("E13").Value = '  +0.86%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.14'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.67%  '
$ws.Range("E15").Value = '  -0.38%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.23'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.44%  '
$ws.Range("D17").Value = '26.575.39'
$ws.Range("E17").Value = '  -0.23%  '
$ws.Range("D18").Value = '0.0₃0743'
$ws.Range("E18").Value = '  +0.55%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '215.82'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.14%  '
$ws.Range("E20").Value = '  +0.21%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.27'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.48%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.35'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.64%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.23'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +14.77%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.33'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.16%  '
$ws.Range("E26").Value = '  +0.20%  '
$ws.Range("E27").Value = '  -0.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.92'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.63%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.62'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.82%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0513'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.20%  '
$ws.Range("E31").Value = '  -0.20%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.36'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.09%  '
$ws.Range("E33").Value = '  +0.91%  '
$ws.Range("D34").Value = '1.259.75'
$ws.Range("E34").Value = '  +7.83%  '
$ws.Range("E35").Value = '  +0.74%  '
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("E37").Value = '  +4.39%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.511'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.44%  '
$ws.Range("B39").Value = 'PaxDollar'
$ws.Range("C39").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.21%  '
$ws.Range("B40").Value = 'ARBITRUM'
$ws.Range("C40").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.798'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.28%  '
$ws.Range("E41").Value = '  -1.66%  '
$ws.Range("E42").Value = '  +0.76%  '
$ws.Range("E43").Value = '  -0.43%  '
$ws.Range("D44").Value = '1.771.72'
$ws.Range("E44").Value = '  +0.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '93.35'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.47%  '
$ws.Range("E46").Value = '  +3.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '55.08'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.77%  '
$ws.Range("D48").Value = '0.0₆0104'
$ws.Range("E48").Value = '  -1.74%  '
$ws.Range("E49").Value = '  +0.28%  '
$ws.Range("E50").Value = '  +0.30%  '
$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.407'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.48%  '
